$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "'19"
$ws.Range("B15").Value = "FR_OPERATIONS"
$ws.Range("C15").Value = "open"
$ws.Range("D15").Value = "2025-03-26T06:37:01Z"
$ws.Range("E15").Value = "bug"
